$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.569.29'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '3.695.89'
$ws.Range("E3").Value = '  +0.44%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '676.69'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.06%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '161.28'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.74%  '
$ws.Range("E9").Value = '  +1.41%  '
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("E11").Value = '  +1.54%  '
$ws.Range("E12").Value = '  +1.01%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '32.57'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("D14").Value = '3.702.81'
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("D15").Value = '69.549.87'
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("E16").Value = '  +2.40%  '
$ws.Range("E17").Value = '  +0.93%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '6.47'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.44%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '471.32'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.80%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '9.85'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.93%  '
$ws.Range("E21").Value = '  +0.98%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '80.53'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.43%  '
$ws.Range("D23").Value = '3.843.62'
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("E24").Value = '  -0.05%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.0000126'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.56%  '
$ws.Range("E26").Value = '  -0.18%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.13'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.23%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.70'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.47%  '
$ws.Range("E29").Value = '  +1.93%  '
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("E32").Value = '  +0.11%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '26.99'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("D34").Value = '3.687.34'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.163'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.44%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '8.45'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +3.75%  '
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("E38").Value = '  -0.01%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("E41").Value = '  +0.43%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '167.23'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("E43").Value = '  +0.30%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '46.55'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.43%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.75'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.51%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.000279'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.16%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '28.07'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("E48").Value = '  -0.17%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.09'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.03%  '
$ws.Range("E50").Value = '  +1.32%  '
$ws.Range("E51").Value = '  +1.91%  '
